# Update the "Max" value for the "pie_threshold_range" row (C5): 28 -> 29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 29

# Move the active selection from C5 to C4, matching the saved sheetView state
$ws.Range("C4").Select()
